$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.967.45'
$ws.Range("E2").Value = '  -1.05%  '

$ws.Range("D3").Value = '1.653.75'
$ws.Range("E3").Value = '  +0.49%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.10'
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3881'
$ws.Range("E7").Value = '  -0.95%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3829'
$ws.Range("E8").Value = '  -1.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.87'
$ws.Range("E9").Value = '  +1.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.354'
$ws.Range("E10").Value = '  -2.70%  '

$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08461'
$ws.Range("E12").Value = '  -1.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.90'
$ws.Range("E13").Value = '  -1.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.107'
$ws.Range("E14").Value = '  -1.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.955'
$ws.Range("E15").Value = '  +2.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001311'
$ws.Range("E16").Value = '  -1.99%  '

$ws.Range("D17").Value = '1.658.76'
$ws.Range("E17").Value = '  +1.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.81'
$ws.Range("E18").Value = '  -0.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06990'
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.76'
$ws.Range("E20").Value = '  -3.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.983'
$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("E23").Value = '  +0.99%  '

$ws.Range("D24").Value = '23.994.85'
$ws.Range("E24").Value = '  -0.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.458'
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.993'
$ws.Range("E26").Value = '  +2.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.17'
$ws.Range("E27").Value = '  -1.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.21'
$ws.Range("E28").Value = '  -3.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.442'
$ws.Range("E29").Value = '  +1.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '138.93'
$ws.Range("E30").Value = '  -1.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.987'
$ws.Range("E31").Value = '  -1.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.530'
$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("D33").Value = '1.839.74'
$ws.Range("E33").Value = '  +1.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.049'
$ws.Range("E34").Value = '  +2.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08107'
$ws.Range("E35").Value = '  -1.41%  '

$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.802'
$ws.Range("E36").Value = '  +0.32%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02946'
$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("E38").Value = '  +0.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2679'
$ws.Range("E39").Value = '  -1.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09131'
$ws.Range("E40").Value = '  -1.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7600'
$ws.Range("E41").Value = '  -1.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.51'
$ws.Range("E42").Value = '  -2.82%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.426'
$ws.Range("E43").Value = '  -1.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.39'
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6996'
$ws.Range("E45").Value = '  -1.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.471'
$ws.Range("E46").Value = '  -1.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.098'
$ws.Range("E47").Value = '  -0.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9997'
$ws.Range("E48").Value = '  +0.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08322'
$ws.Range("E49").Value = '  -0.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.86'
$ws.Range("E50").Value = '  -0.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.206'
$ws.Range("E51").Value = '  -3.26%  '
